# Updated cryptos list (GitHub Actions refresh).
# All Price/Volume cells are plain text in the source workbook (periods used
# as thousands separators, e.g. "27.234.68", plus percentages like
# "  +0.94%  "). For Price values that would otherwise be auto-parsed as a
# number by Excel (dropping trailing zeros / changing precision), a leading
# apostrophe forces the cell to stay text, matching the original inline
# string content exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.234.68'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '1.687.62'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''215.91'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').Value = '''0.520'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').Value = '''22.10'
$ws.Range('E8').Value = '  +8.79%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '''0.259'
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').Value = '''0.0891'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '1.924.35'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '1.692.82'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('D14').Value = '''4.18'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').Value = '''0.553'
$ws.Range('E15').Value = '  +4.73%  '
$ws.Range('D16').Value = '''66.77'
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').Value = '''240.45'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '27.217.97'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').Value = '''8.16'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '0.0₃0746'
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D22').Value = '''4.58'
$ws.Range('E22').Value = '  +3.00%  '
$ws.Range('D23').Value = '''9.57'
$ws.Range('E23').Value = '  +4.00%  '
$ws.Range('E24').Value = '  -3.16%  '
$ws.Range('D25').Value = '''148.28'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').Value = '''7.30'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '''16.45'
$ws.Range('E27').Value = '  +2.19%  '
$ws.Range('D28').Value = '''0.114'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '''0.0501'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('D32').Value = '1.582.01'
$ws.Range('E32').Value = '  +6.74%  '
$ws.Range('D33').Value = '''3.39'
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('D34').Value = '''3.25'
$ws.Range('E34').Value = '  +2.89%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = '''0.963'
$ws.Range('E36').Value = '  +6.77%  '
$ws.Range('D37').Value = '''0.604'
$ws.Range('E37').Value = '  +3.13%  '
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').Value = '''69.66'
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('E43').Value = '  -4.05%  '
$ws.Range('E44').Value = '  -2.56%  '
$ws.Range('D45').Value = '1.832.63'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '''0.788'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').Value = '''91.10'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').Value = '''1.61'
$ws.Range('E48').Value = '  +5.17%  '
$ws.Range('D49').Value = '0.0₆0108'
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('D50').Value = '''0.105'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('E51').Value = '  +5.57%  '
